$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# summer 24 week 12 inputs
$ws.Range("C2").Value = 1.29
$ws.Range("B3").Value = 1.55
$ws.Range("E3").Value = 1.29
$ws.Range("C5").Value = 1.39
$ws.Range("D5").Value = 1.34
$ws.Range("F5").Value = 1.06
